$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in grades (value 5) for newly scored cells
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 5

$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 5

$ws.Range("E28").Value = 5
$ws.Range("F28").Value = 5

# Update the active selection to G28
$ws.Range("G28").Select()
